# Update "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I) columns on the
# Training Dashboard sheet to reflect progress as of 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 29; $row++) {
    $periodCell = $ws.Cells.Item($row, 8)   # column H
    $updateCell = $ws.Cells.Item($row, 9)   # column I

    $periodCell.Value2 = $periodCell.Value2 - 1

    # Keep this a plain text cell (matches the original inline string
    # "04-Nov-2025") instead of letting Excel reinterpret the text as a
    # date serial number. A leading apostrophe forces text entry while
    # preserving the original General number format.
    $updateCell.Value2 = "'04-Nov-2025"
}
